$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E (newest two quarters of data), shifting old D:K to F:M
$ws.Columns("D:E").Insert()

# New D:E columns inherit formats from column F (the former column D) across the used range
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows with no quarterly data (section headers / spacer rows) should not pick up the pasted format
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Populate the two new quarters worth of data, and correct a handful of restated prior-quarter figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2622000
$ws.Range("E8").Value = 2837000
$ws.Range("H8").Value = 2643000
$ws.Range("I8").Value = 2693000
$ws.Range("D9").Value = 1976000
$ws.Range("E9").Value = 2166000
$ws.Range("H9").Value = 1999000
$ws.Range("I9").Value = 2053000
$ws.Range("D10").Value = 646000
$ws.Range("E10").Value = 671000
$ws.Range("H10").Value = 644000
$ws.Range("I10").Value = 640000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 11000
$ws.Range("E14").Value = 85000
$ws.Range("H14").Value = 287000
$ws.Range("I14").Value = 82000
$ws.Range("J14").Value = 106000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 2129000
$ws.Range("E17").Value = 2505000
$ws.Range("H17").Value = 2559000
$ws.Range("I17").Value = 2389000
$ws.Range("D18").Value = 493000
$ws.Range("E18").Value = 332000
$ws.Range("H18").Value = 84000
$ws.Range("I18").Value = 304000
$ws.Range("D20").Value = -147000
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 579000
$ws.Range("E21").Value = 590000
$ws.Range("H21").Value = 369000
$ws.Range("I21").Value = 607000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 346000
$ws.Range("E23").Value = 332000
$ws.Range("H23").Value = 84000
$ws.Range("I23").Value = 304000
$ws.Range("D24").Value = 161000
$ws.Range("E24").Value = 146000
$ws.Range("H24").Value = 705000
$ws.Range("I24").Value = 93000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 185000
$ws.Range("E26").Value = 186000
$ws.Range("H26").Value = -621000
$ws.Range("I26").Value = 211000
$ws.Range("D27").Value = 140000
$ws.Range("E27").Value = 102000
$ws.Range("H27").Value = -644000
$ws.Range("I27").Value = 147000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -12000
$ws.Range("E29").Value = -1000
$ws.Range("H29").Value = -698000
$ws.Range("I29").Value = 5000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 147000
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 128000
$ws.Range("E33").Value = 101000
$ws.Range("H33").Value = -1342000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 128000
$ws.Range("E35").Value = 101000
$ws.Range("H35").Value = -1342000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1166000
$ws.Range("E41").Value = 1187000
$ws.Range("D42").Value = 313000
$ws.Range("E42").Value = 401000
$ws.Range("D43").Value = 1595000
$ws.Range("E43").Value = 1510000
$ws.Range("D44").Value = 577000
$ws.Range("E44").Value = 562000
$ws.Range("D45").Value = 1364000
$ws.Range("E45").Value = 1355000
$ws.Range("D46").Value = 5015000
$ws.Range("E46").Value = 5015000
$ws.Range("D47").Value = 2537000
$ws.Range("E47").Value = 2718000
$ws.Range("D48").Value = 21396000
$ws.Range("E48").Value = 21108000
$ws.Range("D49").Value = 1495000
$ws.Range("E49").Value = 1459000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 2078000
$ws.Range("E52").Value = 2189000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 32521000
$ws.Range("E54").Value = 32489000
$ws.Range("D57").Value = 1329000
$ws.Range("E57").Value = 1299000
$ws.Range("D58").Value = 1659000
$ws.Range("E58").Value = 1308000
$ws.Range("D59").Value = 1411000
$ws.Range("E59").Value = 1440000
$ws.Range("D60").Value = 4399000
$ws.Range("E60").Value = 4047000
$ws.Range("D61").Value = 17636000
$ws.Range("E61").Value = 18088000
$ws.Range("D62").Value = 4003000
$ws.Range("E62").Value = 3766000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 29313000
$ws.Range("E66").Value = 29184000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -1005000
$ws.Range("E72").Value = -1133000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3208000
$ws.Range("E76").Value = 3305000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 128000
$ws.Range("E81").Value = 101000
$ws.Range("H81").Value = -1342000
$ws.Range("D83").Value = 233000
$ws.Range("E83").Value = 258000
$ws.Range("H83").Value = 285000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 680000
$ws.Range("E89").Value = 754000
$ws.Range("H89").Value = 720000
$ws.Range("D91").Value = -529000
$ws.Range("E91").Value = -598000
$ws.Range("H91").Value = -590000
$ws.Range("I91").Value = -464000
$ws.Range("J91").Value = -649000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -315000
$ws.Range("E94").Value = -310000
$ws.Range("H94").Value = -644000
$ws.Range("I94").Value = -859000
$ws.Range("D96").Value = -86000
$ws.Range("E96").Value = -86000
$ws.Range("H96").Value = -79000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -480000
$ws.Range("E100").Value = -434000
$ws.Range("H100").Value = -635000
$ws.Range("D101").Value = -4000
$ws.Range("E101").Value = -30000
$ws.Range("H101").Value = -13000
$ws.Range("I101").Value = 15000
$ws.Range("D102").Value = -119000
$ws.Range("E102").Value = -20000
$ws.Range("H102").Value = -510000
$ws.Range("I102").Value = 417000

